$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after Sheet1 so sheet ordering matches
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Populate header row with strings
$ws2.Range("A1").Value = "four"
$ws2.Range("B1").Value = "five"
$ws2.Range("C1").Value = "six"

# Populate data row with numbers
$ws2.Range("A2").Value = 4
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 6

# Match page margins (inches -> points: 0.75" = 54pt, 1" = 72pt, 0.5" = 36pt)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Select C2 on Sheet2, and make Sheet2 the active/selected tab
$ws2.Range("C2").Select()
$ws2.Activate()
$excel.ActiveWindow.DisplayRuler = $false
